$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51; existing rows 51-86 shift down to 52-87.
$ws.Rows.Item(51).Insert()

# Populate the newly inserted row 51 with the new Maracuya record.
$ws.Cells.Item(51, 1).Value = 10
$ws.Cells.Item(51, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(51, 3).Value = 'La Araucanía'
$ws.Cells.Item(51, 4).Value = 45049
$ws.Cells.Item(51, 5).Value = 9
$ws.Cells.Item(51, 6).Value = 'Fruta'
$ws.Cells.Item(51, 7).Value = 100108
$ws.Cells.Item(51, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(51, 9).Value = 100108003
$ws.Cells.Item(51, 10).Value = 'Maracuyá'
$ws.Cells.Item(51, 11).Value = 'Sin especificar'
$ws.Cells.Item(51, 12).Value = 'Primera'
$ws.Cells.Item(51, 13).Value = 30
$ws.Cells.Item(51, 14).Value = 50000
$ws.Cells.Item(51, 15).Value = 50000
$ws.Cells.Item(51, 16).Value = 50000
$ws.Cells.Item(51, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(51, 18).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(51, 19).Value = 2778
$ws.Cells.Item(51, 20).Value = 18

# Make sure the D51 cell keeps the date number format used by the rest of
# column D (style index 2 in the original workbook).
$ws.Cells.Item(51, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
